$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook window size (cosmetic; best-effort) ---
try {
    $excel.ActiveWindow.WindowWidth = 21570
    $excel.ActiveWindow.WindowHeight = 8055
} catch {
}

# --- Header row: rename "Debit/Credit - Year to date" columns ---
$ws.Range("D1").Value = "Year to date"
$ws.Range("E1").Value = "Past Year"

# --- Column E width ---
$ws.Columns.Item(5).ColumnWidth = 28.74

# --- Data rows: collapse separate Debit/Credit columns into one signed
#     "Year to date" value, repeated into the new "Past Year" column,
#     matching column D's number format on each row ---
$values = @{
    2  = -7875215
    3  = 5415788
    4  = 338000
    5  = 33120
    6  = 20339
    7  = 306921
    8  = 3501263.87
    9  = 1735000
    10 = 553539
    11 = 471000
    12 = 15238000
    13 = -235000
    14 = -3500
    15 = -5718554
    16 = -13302057
    17 = -478644.87
}

foreach ($r in 2..17) {
    $v = $values[$r]
    $ws.Range("D$r").Value = $v
    $ws.Range("E$r").Value = $v
    $ws.Range("E$r").NumberFormat = $ws.Range("D$r").NumberFormat
}

# --- Selection as left by the editor ---
$ws.Range("I17").Select()
